$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text is a valid numeric literal (e.g. "0.07660") need the
# cell pre-formatted as Text, otherwise Excel auto-converts the input to a number
# and silently drops significant trailing zeros (e.g. "0.07660" -> 0.0766).
$ws.Range('D2').Value = '27.450.00'
$ws.Range('E2').Value = '  +5.00%  '
$ws.Range('D3').Value = '1.722.58'
$ws.Range('E3').Value = '  +4.40%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.28'
$ws.Range('E5').Value = '  +3.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5336'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2653'
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06592'
$ws.Range('E9').Value = '  +4.50%  '
$ws.Range('E10').Value = '  +5.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07660'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.595'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '1.737.22'
$ws.Range('E13').Value = '  +5.20%  '
$ws.Range('D14').Value = '1.960.08'
$ws.Range('E14').Value = '  +4.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5786'
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('D16').Value = '0.0₅8270'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.68'
$ws.Range('E17').Value = '  +3.88%  '
$ws.Range('D18').Value = '27.459.07'
$ws.Range('E18').Value = '  +5.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.21'
$ws.Range('E19').Value = '  +11.75%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.709'
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.55'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.009'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.03'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.734'
$ws.Range('E26').Value = '  +13.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1230'
$ws.Range('E27').Value = '  +3.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.316'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.45'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05467'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.298'
$ws.Range('E31').Value = '  +2.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.542'
$ws.Range('E32').Value = '  +3.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.426'
$ws.Range('E33').Value = '  +2.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.658'
$ws.Range('E34').Value = '  +6.15%  '
$ws.Range('E35').Value = '  +2.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9550'
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5908'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01641'
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.903'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.047.14'
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8448'
$ws.Range('E42').Value = '  +3.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.18'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '1.867.91'
$ws.Range('E45').Value = '  +4.55%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  +9.20%  '
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4500'
$ws.Range('E48').Value = '  +4.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.169'
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.003'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  +2.66%  '
